$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.185.13"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.322.20"
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "302.81"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "99.44"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.19%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.516"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.69%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "36.08"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.92%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  -1.15%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "17.55"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.75%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.91"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("D15").Value = "2.683.28"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "2.310.45"
$ws.Range("E16").Value = "  +1.26%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.797"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "43.115.32"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.21"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +7.04%  "
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "0.0₃0911"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "68.13"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.49%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "241.69"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("E26").Value = "  -0.17%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "25.46"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.05%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "168.29"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  +1.44%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.20"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("E32").Value = "  +3.74%  "
$ws.Range("E33").Value = "  -0.07%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.73"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.39%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "17.85"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +5.38%  "
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").Value = "1.999.84"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("E44").Value = "  -4.91%  "
$ws.Range("E45").Value = "  +0.49%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "17.65"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("E47").Value = "  +0.18%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "76.64"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +9.09%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "54.98"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.87"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +12.97%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.548.07"
$ws.Range("E51").Value = "  +0.82%  "
